$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (B7, C7, E7) with the new "/if" route info
$ws.Range("B7").Value = "/if"
$ws.Range("C7").Value = "ifexist"
$ws.Range("E7").Value = "check if a symbol exist"

# Update selection to I18
$ws.Range("I18").Select()
